# Applies the cryptos-list refresh described in the commit:
# "Updated cryptos list on Thu Jan  4 20:36:59 UTC 2024 with GitHub Actions"
#
# Strategy: write each changed cell's new text verbatim via Range.Value.
# Column D sometimes holds numeric-looking text (e.g. "320.14") that must
# stay a text string (matches the source inlineStr/shared-string cells,
# not numbers) - exactly like the workbook already stores "44.145.44" or
# "0.591" as text. A leading apostrophe forces Excel to keep such values
# as text instead of auto-converting them to numbers, without touching
# any other formatting.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '44.244.49'
$ws.Range("E2").Value = '  +3.04%  '
# Row 3
$ws.Range("D3").Value = '2.283.54'
$ws.Range("E3").Value = '  +2.92%  '
# Row 4
$ws.Range("E4").Value = '  -0.07%  '
# Row 5
$ws.Range("D5").Value = '''320.14'
$ws.Range("E5").Value = '  +1.22%  '
# Row 6
$ws.Range("D6").Value = '''107.54'
$ws.Range("E6").Value = '  +8.10%  '
# Row 7
$ws.Range("D7").Value = '''0.592'
$ws.Range("E7").Value = '  +0.31%  '
# Row 8
$ws.Range("E8").Value = '  -0.14%  '
# Row 9
$ws.Range("D9").Value = '''0.576'
$ws.Range("E9").Value = '  +2.45%  '
# Row 10
$ws.Range("D10").Value = '''39.30'
$ws.Range("E10").Value = '  +6.22%  '
# Row 11
$ws.Range("D11").Value = '''0.0845'
$ws.Range("E11").Value = '  +2.08%  '
# Row 12
$ws.Range("D12").Value = '''7.96'
$ws.Range("E12").Value = '  +2.22%  '
# Row 13
$ws.Range("E13").Value = '  +1.54%  '
# Row 14
$ws.Range("D14").Value = '''0.890'
$ws.Range("E14").Value = '  +3.16%  '
# Row 15
$ws.Range("D15").Value = '2.627.14'
$ws.Range("E15").Value = '  +2.79%  '
# Row 16
$ws.Range("D16").Value = '''14.71'
$ws.Range("E16").Value = '  +3.36%  '
# Row 17
$ws.Range("D17").Value = '2.281.77'
$ws.Range("E17").Value = '  +3.32%  '
# Row 18
$ws.Range("D18").Value = '44.147.34'
# Row 19
$ws.Range("D19").Value = '''14.14'
$ws.Range("E19").Value = '  -8.47%  '
# Row 20
$ws.Range("E20").Value = '  +4.12%  '
# Row 21
$ws.Range("D21").Value = '''6.57'
$ws.Range("E21").Value = '  +1.81%  '
# Row 22
$ws.Range("D22").Value = '''66.59'
$ws.Range("E22").Value = '  +1.79%  '
# Row 23
$ws.Range("D23").Value = '''3.23'
$ws.Range("E23").Value = '  +1.66%  '
# Row 24
$ws.Range("D24").Value = '''239.06'
$ws.Range("E24").Value = '  +0.99%  '
# Row 25
$ws.Range("D25").Value = '''2.22'
$ws.Range("E25").Value = '  +5.08%  '
# Row 26
$ws.Range("E26").Value = '  -0.17%  '
# Row 27
$ws.Range("D27").Value = '''10.28'
$ws.Range("E27").Value = '  +1.82%  '
# Row 28
$ws.Range("D28").Value = '''39.26'
$ws.Range("E28").Value = '  +14.77%  '
# Row 29
$ws.Range("D29").Value = '''2.23'
$ws.Range("E29").Value = '  +0.66%  '
# Row 30
$ws.Range("E30").Value = '  +3.50%  '
# Row 31
$ws.Range("B31").Value = 'Monero'
$ws.Range("C31").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D31").Value = '''163.45'
$ws.Range("E31").Value = '  +3.77%  '
# Row 32
$ws.Range("B32").Value = 'EthereumClassic'
$ws.Range("C32").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D32").Value = '''20.73'
$ws.Range("E32").Value = '  +1.08%  '
# Row 33
$ws.Range("D33").Value = '''0.0892'
$ws.Range("E33").Value = '  +1.18%  '
# Row 34
$ws.Range("E34").Value = '  -1.27%  '
# Row 35
$ws.Range("B35").Value = 'LidoDAOToken'
$ws.Range("C35").Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range("D35").Value = '''3.30'
$ws.Range("E35").Value = '  +2.77%  '
# Row 36
$ws.Range("B36").Value = 'ARBITRUM'
$ws.Range("C36").Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range("D36").Value = '''2.08'
$ws.Range("E36").Value = '  +5.00%  '
# Row 37
$ws.Range("E37").Value = '  +11.93%  '
# Row 38
$ws.Range("D38").Value = '''0.121'
$ws.Range("E38").Value = '  -0.69%  '
# Row 39
$ws.Range("D39").Value = '''4.01'
$ws.Range("E39").Value = '  +6.34%  '
# Row 40
$ws.Range("D40").Value = '''4.50'
$ws.Range("E40").Value = '  +1.10%  '
# Row 41
$ws.Range("E41").Value = '  +1.39%  '
# Row 42
$ws.Range("D42").Value = '''15.39'
$ws.Range("E42").Value = '  +24.77%  '
# Row 43
$ws.Range("E43").Value = '  +0.09%  '
# Row 44
$ws.Range("D44").Value = '1.773.06'
$ws.Range("E44").Value = '  -7.18%  '
# Row 45
$ws.Range("E45").Value = '  +0.94%  '
# Row 46
$ws.Range("D46").Value = '''86.21'
$ws.Range("E46").Value = '  -3.53%  '
# Row 47
$ws.Range("D47").Value = '''5.45'
$ws.Range("E47").Value = '  +0.63%  '
# Row 48
$ws.Range("D48").Value = '''75.69'
$ws.Range("E48").Value = '  +0.07%  '
# Row 49
$ws.Range("D49").Value = '''60.16'
$ws.Range("E49").Value = '  -1.10%  '
# Row 50
$ws.Range("E50").Value = '  +2.42%  '
# Row 51
$ws.Range("E51").Value = '  +6.34%  '
